# Update Jag1-Notch1 LR-pair metrics following Dr Hou advice
# (recomputed with 3 ligand/receptor-expressing-cell replicates instead of 1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 16.12468433333333
$ws.Range("H2").Value = 48.374053
$ws.Range("I2").Value = 0.2955490655206278
$ws.Range("J2").Value = 0.2955490655206279
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 68.00339
$ws.Range("N2").Value = 204.01017
$ws.Range("O2").Value = 0.6265962299909886
$ws.Range("P2").Value = 0.6265962299909885
$ws.Range("Q2").Value = 1096.533197346556
$ws.Range("R2").Value = 9868.798776119009
$ws.Range("S2").Value = 0.1851899302325851
$ws.Range("T2").Value = 0.1851899302325851
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 16.12468433333333
$ws.Range("H3").Value = 48.374053
$ws.Range("I3").Value = 0.2955490655206278
$ws.Range("J3").Value = 0.2955490655206279
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.991529999999999
$ws.Range("N3").Value = 26.97459
$ws.Range("O3").Value = 0.08284967558015671
$ws.Range("P3").Value = 0.08284967558015671
$ws.Range("Q3").Value = 144.9855829236966
$ws.Range("R3").Value = 1304.87024631327
$ws.Range("S3").Value = 0.0244861441964025
$ws.Range("T3").Value = 0.0244861441964025
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 16.12468433333333
$ws.Range("H4").Value = 48.374053
$ws.Range("I4").Value = 0.2955490655206278
$ws.Range("J4").Value = 0.2955490655206279
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.76843933333333
$ws.Range("N4").Value = 32.305318
$ws.Range("O4").Value = 0.09922245772090688
$ws.Range("P4").Value = 0.09922245772090688
$ws.Range("Q4").Value = 173.6376850126504
$ws.Range("R4").Value = 1562.739165113854
$ws.Range("S4").Value = 0.02932510465807403
$ws.Range("T4").Value = 0.02932510465807403
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.12468433333333
$ws.Range("H5").Value = 48.374053
$ws.Range("I5").Value = 0.2955490655206278
$ws.Range("J5").Value = 0.2955490655206279
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 20.764887
$ws.Range("N5").Value = 62.294661
$ws.Range("O5").Value = 0.1913316367079478
$ws.Range("P5").Value = 0.1913316367079478
$ws.Range("Q5").Value = 334.8272480923369
$ws.Range("R5").Value = 3013.445232831033
$ws.Range("S5").Value = 0.05654788643356622
$ws.Range("T5").Value = 0.05654788643356623
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 14.68975
$ws.Range("H6").Value = 44.06925
$ws.Range("I6").Value = 0.2692481784748309
$ws.Range("J6").Value = 0.2692481784748309
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 68.00339
$ws.Range("N6").Value = 204.01017
$ws.Range("O6").Value = 0.6265962299909886
$ws.Range("P6").Value = 0.6265962299909885
$ws.Range("Q6").Value = 998.9527982524999
$ws.Range("R6").Value = 8990.575184272499
$ws.Range("S6").Value = 0.1687098935642699
$ws.Range("T6").Value = 0.1687098935642699
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 14.68975
$ws.Range("H7").Value = 44.06925
$ws.Range("I7").Value = 0.2692481784748309
$ws.Range("J7").Value = 0.2692481784748309
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.991529999999999
$ws.Range("N7").Value = 26.97459
$ws.Range("O7").Value = 0.08284967558015671
$ws.Range("P7").Value = 0.08284967558015671
$ws.Range("Q7").Value = 132.0833278175
$ws.Range("R7").Value = 1188.7499503575
$ws.Range("S7").Value = 0.02230712423718787
$ws.Range("T7").Value = 0.02230712423718787
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.68975
$ws.Range("H8").Value = 44.06925
$ws.Range("I8").Value = 0.2692481784748309
$ws.Range("J8").Value = 0.2692481784748309
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.76843933333333
$ws.Range("N8").Value = 32.305318
$ws.Range("O8").Value = 0.09922245772090688
$ws.Range("P8").Value = 0.09922245772090688
$ws.Range("Q8").Value = 158.1856816968333
$ws.Range("R8").Value = 1423.6711352715
$ws.Range("S8").Value = 0.0267154660051501
$ws.Range("T8").Value = 0.0267154660051501
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.68975
$ws.Range("H9").Value = 44.06925
$ws.Range("I9").Value = 0.2692481784748309
$ws.Range("J9").Value = 0.2692481784748309
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 20.764887
$ws.Range("N9").Value = 62.294661
$ws.Range("O9").Value = 0.1913316367079478
$ws.Range("P9").Value = 0.1913316367079478
$ws.Range("Q9").Value = 305.0309988082499
$ws.Range("R9").Value = 2745.27898927425
$ws.Range("S9").Value = 0.05151569466822303
$ws.Range("T9").Value = 0.05151569466822303
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.967860666666667
$ws.Range("H10").Value = 5.903582
$ws.Range("I10").Value = 0.03606888476606249
$ws.Range("J10").Value = 0.03606888476606249
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 68.00339
$ws.Range("N10").Value = 204.01017
$ws.Range("O10").Value = 0.6265962299909886
$ws.Range("P10").Value = 0.6265962299909885
$ws.Range("Q10").Value = 133.8211963809933
$ws.Range("R10").Value = 1204.39076742894
$ws.Range("S10").Value = 0.02260062721439416
$ws.Range("T10").Value = 0.02260062721439415
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.967860666666667
$ws.Range("H11").Value = 5.903582
$ws.Range("I11").Value = 0.03606888476606249
$ws.Range("J11").Value = 0.03606888476606249
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.991529999999999
$ws.Range("N11").Value = 26.97459
$ws.Range("O11").Value = 0.08284967558015671
$ws.Range("P11").Value = 0.08284967558015671
$ws.Range("Q11").Value = 17.69407822015333
$ws.Range("R11").Value = 159.24670398138
$ws.Range("S11").Value = 0.002988295401406334
$ws.Range("T11").Value = 0.002988295401406334
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.967860666666667
$ws.Range("H12").Value = 5.903582
$ws.Range("I12").Value = 0.03606888476606249
$ws.Range("J12").Value = 0.03606888476606249
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.76843933333333
$ws.Range("N12").Value = 32.305318
$ws.Range("O12").Value = 0.09922245772090688
$ws.Range("P12").Value = 0.09922245772090688
$ws.Range("Q12").Value = 21.19078820545289
$ws.Range("R12").Value = 190.717093849076
$ws.Range("S12").Value = 0.003578843393740897
$ws.Range("T12").Value = 0.003578843393740897
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.967860666666667
$ws.Range("H13").Value = 5.903582
$ws.Range("I13").Value = 0.03606888476606249
$ws.Range("J13").Value = 0.03606888476606249
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 20.764887
$ws.Range("N13").Value = 62.294661
$ws.Range("O13").Value = 0.1913316367079478
$ws.Range("P13").Value = 0.1913316367079478
$ws.Range("Q13").Value = 40.862404375078
$ws.Range("R13").Value = 367.761639375702
$ws.Range("S13").Value = 0.0069011187565211
$ws.Range("T13").Value = 0.0069011187565211
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 21.77610566666667
$ws.Range("H14").Value = 65.328317
$ws.Range("I14").Value = 0.3991338712384788
$ws.Range("J14").Value = 0.3991338712384788
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 68.00339
$ws.Range("N14").Value = 204.01017
$ws.Range("O14").Value = 0.6265962299909886
$ws.Range("P14").Value = 0.6265962299909885
$ws.Range("Q14").Value = 1480.849006331543
$ws.Range("R14").Value = 13327.64105698389
$ws.Range("S14").Value = 0.2500957789797395
$ws.Range("T14").Value = 0.2500957789797394
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 21.77610566666667
$ws.Range("H15").Value = 65.328317
$ws.Range("I15").Value = 0.3991338712384788
$ws.Range("J15").Value = 0.3991338712384788
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 8.991529999999999
$ws.Range("N15").Value = 26.97459
$ws.Range("O15").Value = 0.08284967558015671
$ws.Range("P15").Value = 0.08284967558015671
$ws.Range("Q15").Value = 195.8005073850033
$ws.Range("R15").Value = 1762.20456646503
$ws.Range("S15").Value = 0.03306811174516001
$ws.Range("T15").Value = 0.03306811174516001
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 21.77610566666667
$ws.Range("H16").Value = 65.328317
$ws.Range("I16").Value = 0.3991338712384788
$ws.Range("J16").Value = 0.3991338712384788
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 10.76843933333333
$ws.Range("N16").Value = 32.305318
$ws.Range("O16").Value = 0.09922245772090688
$ws.Range("P16").Value = 0.09922245772090688
$ws.Range("Q16").Value = 234.4946727877562
$ws.Range("R16").Value = 2110.452055089806
$ws.Range("S16").Value = 0.03960304366394185
$ws.Range("T16").Value = 0.03960304366394185
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 21.77610566666667
$ws.Range("H17").Value = 65.328317
$ws.Range("I17").Value = 0.3991338712384788
$ws.Range("J17").Value = 0.3991338712384788
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 20.764887
$ws.Range("N17").Value = 62.294661
$ws.Range("O17").Value = 0.1913316367079478
$ws.Range("P17").Value = 0.1913316367079478
$ws.Range("Q17").Value = 452.1783734683929
$ws.Range("R17").Value = 4069.605361215537
$ws.Range("S17").Value = 0.07636693684963743
$ws.Range("T17").Value = 0.07636693684963743
